$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 9575
$ws1.Range("F10").Value = 2208
$ws1.Range("F14").Value = 273
$ws1.Range("F18").Value = 1329

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 16

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 16
$ws4.Range("F7").Value = 9575
$ws4.Range("F11").Value = 2208
$ws4.Range("F15").Value = 273
$ws4.Range("F19").Value = 1329
